# ---------------------------------------------------------------------------
# Edit: simplify/restructure GST report -> "GST Audit Report"
# Rebuilds headers, data row, TOTAL row, and footer notes with the
# Arial-based styling (header band, thin-bordered table, bold total row).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename sheet -----------------------------------------------------------
$ws.Name = "GST Audit Report"

# offset Excel adds internally to ColumnWidth when persisting to OOXML width
$colOffset = 0.8333333333333334

# --- column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 8  - $colOffset
$ws.Columns.Item(2).ColumnWidth  = 25 - $colOffset
$ws.Columns.Item(3).ColumnWidth  = 18 - $colOffset
$ws.Columns.Item(4).ColumnWidth  = 20 - $colOffset
$ws.Columns.Item(5).ColumnWidth  = 12 - $colOffset
$ws.Columns.Item(6).ColumnWidth  = 15 - $colOffset
$ws.Columns.Item(7).ColumnWidth  = 12 - $colOffset
$ws.Columns.Item(8).ColumnWidth  = 12 - $colOffset
$ws.Columns.Item(9).ColumnWidth  = 12 - $colOffset
$ws.Columns.Item(10).ColumnWidth = 12 - $colOffset
$ws.Columns.Item(11).ColumnWidth = 35 - $colOffset

# --- clear any previous content/formatting on the used range ---------------
$ws.Cells.Clear()

# --- header row (row 1) ------------------------------------------------------
$headers = @("S.No.", "Vendor/Shop Name", "GSTIN", "Invoice No.", "Date", "Taxable Amount", "Total Tax", "CGST", "SGST", "IGST", "HSN Codes")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 11
    $cell.Font.Bold = $true
    $cell.Font.Color = 16777215
    $cell.Interior.Color = 8674603
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}
$ws.Rows.Item(1).RowHeight = 30

# --- data row (row 2) --------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 25

$a2 = $ws.Cells.Item(2, 1)
$a2.Value = 1
$a2.HorizontalAlignment = -4108

$b2 = $ws.Cells.Item(2, 2)
$b2.Value = "WESTSIDE, Sjr Zion, Survey"
$b2.Font.Name = "Arial"
$b2.Font.Size = 10
$b2.VerticalAlignment = -4108
$b2.WrapText = $true
$b2.Borders.LineStyle = 1

$dataVals = @("29AAACL1838J1ZC", "W089 100169940", "2024-09-28", "4045.01", "173.91", "173.91", "173.91", "N/A")
# these columns hold numeric/date-looking text that must stay plain text
$dataForceText = @($false, $false, $true, $true, $true, $true, $true, $false)
for ($i = 0; $i -lt $dataVals.Length; $i++) {
    $cell = $ws.Cells.Item(2, $i + 3)
    if ($dataForceText[$i]) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $dataVals[$i]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.VerticalAlignment = -4108
    $cell.Borders.LineStyle = 1
}

$k2 = $ws.Cells.Item(2, 11)
$k2.Value = "996211, 62052000, 62052000, 62046200, 48194000, 33072000, 39264099"
$k2.Font.Name = "Arial"
$k2.Font.Size = 10
$k2.VerticalAlignment = -4108
$k2.WrapText = $true
$k2.Borders.LineStyle = 1

# --- TOTAL row (row 3) --------------------------------------------------------
$a3 = $ws.Cells.Item(3, 1)
$a3.Value = "TOTAL"
$a3.Borders.LineStyle = 1

# top+bottom borders on the blank merged-in cells (left as default font,
# they stay empty so only the border is ever visible)
$ws.Cells.Item(3, 2).Borders.Item(8).LineStyle = 1
$ws.Cells.Item(3, 2).Borders.Item(9).LineStyle = 1
$ws.Cells.Item(3, 3).Borders.Item(8).LineStyle = 1
$ws.Cells.Item(3, 3).Borders.Item(9).LineStyle = 1
$ws.Cells.Item(3, 4).Borders.Item(8).LineStyle = 1
$ws.Cells.Item(3, 4).Borders.Item(9).LineStyle = 1

# top+bottom+right border on the right-most cell of the merge
$e3 = $ws.Cells.Item(3, 5)
$e3.Borders.Item(8).LineStyle = 1
$e3.Borders.Item(9).LineStyle = 1
$e3.Borders.Item(10).LineStyle = 1

# merge AFTER the per-cell border edges are set, matching Excel's own
# behaviour of collapsing internal shared edges once a range is merged
$ws.Range("A3:E3").Merge()

# re-assert A3's right edge: merging silently drops it because the
# neighbouring (now absorbed) cell did not carry a right border of its own
$a3.Borders.Item(10).LineStyle = 1

# apply A3's font/alignment AFTER merging: the merge operation otherwise
# propagates the anchor cell's font onto the absorbed cells as well
$a3.Font.Name = "Arial"
$a3.Font.Size = 10
$a3.Font.Bold = $true
$a3.HorizontalAlignment = -4108

$totalVals = @("₹4,045.01", "₹173.91", "₹173.91", "₹173.91", "₹0.00")
for ($i = 0; $i -lt $totalVals.Length; $i++) {
    $cell = $ws.Cells.Item(3, $i + 6)
    $cell.NumberFormat = "@"
    $cell.Value = $totalVals[$i]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4152
    $cell.Borders.LineStyle = 1
}

# --- footer notes (rows 5-6) --------------------------------------------------
$a5 = $ws.Cells.Item(5, 1)
$a5.Value = "Generated for GST Audit Purposes"
$a5.Font.Name = "Arial"
$a5.Font.Size = 9
$a5.Font.Italic = $true

$a6 = $ws.Cells.Item(6, 1)
$a6.Value = "Total Invoices: 1"
$a6.Font.Name = "Arial"
$a6.Font.Size = 9

# --- freeze header row, keep selection on A1 ---------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

Write-Host "GST Audit Report rebuilt"
